# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces two changes from the source commit:
#   1. The table on slide 5 gets its table-style swapped from the
#      custom "Table_0" style to the built-in style
#      {857C5A2D-4BEE-4AB5-B34E-2218163C544B}.
#   2. The theme used by the slide master (stored as ppt/theme/theme2.xml,
#      color scheme name "Red Violet"/"Integral") is recoloured to match
#      the standard Office theme palette (the palette that, before the
#      edit, lived in ppt/theme/theme1.xml - the Notes Master's theme).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Table style swap (slide 5, the table graphic frame).
# ---------------------------------------------------------------------
$oldStyleId = "{0F695865-D3F1-4CB3-B0A0-933D15CA9AA3}"
$newStyleId = "{857C5A2D-4BEE-4AB5-B34E-2218163C544B}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colour swap - recolour the active theme (the slide master's
#    theme) from the "Integral" / "Red Violet" palette to the standard
#    "Office" palette.
# ---------------------------------------------------------------------
function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette (the standard "Office" theme colours), in the
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order used by
# ThemeColorScheme.Item(1..12).
$officeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slideForTheme = $p.Slides.Item(1)
$themeColors = $slideForTheme.ThemeColorScheme

for ($i = 1; $i -le $officeHex.Count; $i++) {
    $themeColors.Item($i).RGB = ConvertTo-OleColor $officeHex[$i - 1]
}
